$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '67.816.41'
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '3.541.76'
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = '  -3.62%  '
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.10%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '616.18'
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = '  -5.02%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '154.04'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  -3.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '3.535.73'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  -3.61%  '
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  +0.02%  '
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -2.09%  '
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  -2.37%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.92'
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = '  -1.89%  '
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  -1.61%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000225'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -2.45%  '
$ws.Range("B14").NumberFormat = "@"
$ws.Range("B14").Value = 'Avalanche'
$ws.Range("C14").NumberFormat = "@"
$ws.Range("C14").Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '32.26'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  -0.61%  '
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = 'WrappedliquidstakedEther2.0'
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '4.143.47'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  -3.40%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '3.530.04'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -4.02%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '67.770.02'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  -2.40%  '
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.33%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '6.40'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  -0.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '15.54'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  -2.75%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '454.82'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  -2.30%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '9.42'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  -2.98%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.642'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  +0.09%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '77.71'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  -2.07%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.684.65'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -3.51%  '
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +0.10%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0000119'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  -4.49%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '10.51'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  -2.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '8.38'
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -5.94%  '
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  -1.57%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.64'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  -1.07%  '
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  -0.55%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '25.97'
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -2.26%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.91'
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -4.06%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.25'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  -3.24%  '
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  -2.26%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.542.40'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  -3.28%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '8.02'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -3.90%  '
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.04%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.07%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '176.35'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -1.58%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '5.64'
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = '  -4.59%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0883'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  -1.04%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.08'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  -4.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.889'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  -4.39%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.60'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +9.68%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '45.91'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  -1.84%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.61'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  -3.55%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '7.69'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -1.22%  '
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -3.79%  '
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  -3.65%  '
